$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title/name of the unit in row 3
$ws.Range("B3").Value = "Hladilni agregat Climaveneta NECS SL 1314"

# Update technical data rows 9-21
$ws.Range("B9").Value  = "Hladilna moč:    332,3 kW"
$ws.Range("B10").Value = "EER (EN14511 metoda):    2,55 "
$ws.Range("B11").Value = "ESEER (EN14511 metoda):    4,10 "
$ws.Range("B12").Value = "SEER (Reg. EU 2016/2281):    4,03 "
$ws.Range("B13").Value = "El. priključek:    400V/ 3F/ 50Hz "
$ws.Range("B14").Value = "Zvočni tlak (SPL):    54 dB(A)"
$ws.Range("B15").Value = "Zvočna moč (PWL):    86 dB(A)"
$ws.Range("B16").Value = "Število hladilnih krogov:    2 "
$ws.Range("B17").Value = "Število kompresorjev:    4 "
$ws.Range("B18").Value = "Dolžina:    5080 mm"
$ws.Range("B19").Value = "Širina:    2260 mm"
$ws.Range("B20").Value = "Višina:    2450 mm"
$ws.Range("B21").Value = "Teža:    3060 kg"

# Update quantity / price columns for row 21. Quantity (C21) and unit price
# (D21) are recorded as text (matching the source data), total price (E21)
# as a real number. Briefly flip to a Text format while assigning so Excel
# doesn't auto-convert the numeric-looking strings, then restore each
# cell's original number format.
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "2"
$ws.Range("C21").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1000"
$ws.Range("D21").NumberFormat = "0.00"

$ws.Range("E21").Value = 2000

# Remove the second unit's section (rows 22-42), which no longer exists
$ws.Range("A22:E42").EntireRow.Delete()
